$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header block (row 1) ---
$ws.Range("A1").Value = "Cotización de Productos"
$ws.Range("C1").Value = "Cotizacion N°"
$ws.Range("D1").Value = 60023

# --- 2. Client / provider info updates ---
$ws.Range("B4").Value = "Mundo Patitas SAC"
$ws.Range("B7").Value = 987654321
$ws.Range("B12").Value = "Av. Ejemplo 123, Los Olivos"

# --- 3. Insert 5 new rows before row 20 to make room for the extra product lines ---
$ws.Rows.Item(20).Resize(5).Insert()

# Carry the detail-row formatting (borders/number format) down into the new rows
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E24").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# --- 4. Rewrite the product detail rows (17-24), now sorted A-Z by Descripcion ---
$ws.Range("A17").Value = 1019
$ws.Range("B17").Value = "Cepillo Dental Para Perro Medium"
$ws.Range("C17").Value = 40
$ws.Range("D17").Value = 10
$ws.Range("E17").Formula = "=IFERROR(C17*D17,`"`")"

$ws.Range("A18").Value = 1045
$ws.Range("B18").Value = "Filtro para Acuario 3L"
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 8
$ws.Range("E18").Formula = "=IFERROR(C18*D18,`"`")"

$ws.Range("A19").Value = 1029
$ws.Range("B19").Value = "Snack Mixto Para Aves Tropicales 500g"
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 4
$ws.Range("E19").Formula = "=IFERROR(C19*D19,`"`")"

$ws.Range("A20").Value = 1036
$ws.Range("B20").Value = "Snack Saludable Para Perro 100g"
$ws.Range("C20").Value = 20
$ws.Range("D20").Value = 10.5
$ws.Range("E20").Formula = "=IFERROR(C20*D20,`"`")"

$ws.Range("A21").Value = 1047
$ws.Range("B21").Value = "Suéter para Perro Mediano"
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 5.5
$ws.Range("E21").Formula = "=IFERROR(C21*D21,`"`")"

$ws.Range("A22").Value = 1046
$ws.Range("B22").Value = "Termómetro Digital Acuario"
$ws.Range("C22").Value = 30
$ws.Range("D22").Value = 8
$ws.Range("E22").Formula = "=IFERROR(C22*D22,`"`")"

$ws.Range("A23").Value = 1026
$ws.Range("B23").Value = "Vitaminas Multinivel para Perro 250g"
$ws.Range("C23").Value = 100
$ws.Range("D23").Value = 5
$ws.Range("E23").Formula = "=IFERROR(C23*D23,`"`")"

$ws.Range("A24").Value = 1027
$ws.Range("B24").Value = "Vitaminas de Salmón para Gato 500ml"
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 6.5
$ws.Range("E24").Formula = "=IFERROR(C24*D24,`"`")"

# --- 5. Fix the SubTotal formula to cover the expanded range ---
$ws.Range("E26").Formula = "=SUM(E17:E24)"

# --- 6. Expand the "Detalle" table over the new range ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A16:E24"))

# --- 7. View settings (zoom/top-left cell/selection) ---
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("C20").Select()
